$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A54").Value = 45986
$ws.Range("A53").Copy()
$ws.Range("A54").PasteSpecial(-4122)
$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = -2.451276118722334
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = -0.8888225292121632
